$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.890.63'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.736.94'
$ws.Range('E3').Value = '  +3.73%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.77'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.19'
$ws.Range('E6').Value = '  +6.80%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.548'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.736.02'
$ws.Range('E9').Value = '  +3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.73'
$ws.Range('E14').Value = '  +3.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.235.34'
$ws.Range('E15').Value = '  +3.77%  '
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.827.66'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.735.90'
$ws.Range('E18').Value = '  +4.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.84'
$ws.Range('E19').Value = '  +4.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '373.19'
$ws.Range('E20').Value = '  +3.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.70'
$ws.Range('E21').Value = '  +5.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.52'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.95'
$ws.Range('E23').Value = '  +4.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  +3.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.84'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('E27').Value = '  +2.61%  '
$ws.Range('E28').Value = '  +3.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000107'
$ws.Range('E29').Value = '  +3.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '593.72'
$ws.Range('E30').Value = '  +6.47%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  +4.44%  '
$ws.Range('E33').Value = '  +5.36%  '
$ws.Range('E34').Value = '  +6.60%  '
$ws.Range('E35').Value = '  +5.11%  '
$ws.Range('E36').Value = '  +5.50%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.36'
$ws.Range('E38').Value = '  +1.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.86'
$ws.Range('E39').Value = '  +1.13%  '
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.51'
$ws.Range('E42').Value = '  +4.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.68'
$ws.Range('E43').Value = '  +3.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.99'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0314'
$ws.Range('E46').Value = '  -2.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.92'
$ws.Range('E47').Value = '  +1.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.85'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.96'
$ws.Range('E49').Value = '  +6.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.80'
$ws.Range('E50').Value = '  +7.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.606'
$ws.Range('E51').Value = '  +7.61%  '
